$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G column (bold flag) to FALSE for all footnote rows (column C = "footnote")
$ws.Range("G6:G7").Value = $false
$ws.Range("G12:G14").Value = $false
$ws.Range("G19:G21").Value = $false
$ws.Range("G26:G28").Value = $false
$ws.Range("G33:G36").Value = $false
$ws.Range("G41:G44").Value = $false
$ws.Range("G49:G52").Value = $false
$ws.Range("G57:G60").Value = $false
$ws.Range("G65:G68").Value = $false
$ws.Range("G73:G76").Value = $false
$ws.Range("G81:G85").Value = $false
$ws.Range("G90:G93").Value = $false
$ws.Range("G98:G101").Value = $false
$ws.Range("G106").Value = $false
$ws.Range("G111:G112").Value = $false
$ws.Range("G117:G120").Value = $false
$ws.Range("G125:G126").Value = $false
$ws.Range("G131:G133").Value = $false
$ws.Range("G138:G143").Value = $false
$ws.Range("G148:G153").Value = $false
$ws.Range("G158:G159").Value = $false
$ws.Range("G164:G165").Value = $false
$ws.Range("G170:G171").Value = $false
$ws.Range("G176:G178").Value = $false
$ws.Range("G183:G187").Value = $false
$ws.Range("G192:G196").Value = $false
$ws.Range("G201:G202").Value = $false
$ws.Range("G207:G208").Value = $false
$ws.Range("G213:G214").Value = $false
$ws.Range("G219").Value = $false
$ws.Range("G224:G225").Value = $false

# Apply AutoFilter over the full data range
$ws.Range("A1:I225").AutoFilter()

# Register the hidden _FilterDatabase defined name created by AutoFilter
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$I`$225")
$filterName.Visible = $false

# Update the sheet view: scroll back to top-left and move the selection
$ws.Range("F17").Select()
